$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Robot2 location change (2, 8) -> (11, 8)
$ws.Range("A1").Value = "Move Robot2 to location (11, 8) and remove the toolkit."
$ws.Range("E1").Value = "(11, 8)"

# Row 2: Robot26 location change (11, 4) -> (4, 4)
$ws.Range("A2").Value = "Move Robot26 to location (4, 4) and remove the liquid spill."
$ws.Range("E2").Value = "(4, 4)"

# Row 3: Robot42 location change (9, 5) -> (9, 1)
$ws.Range("A3").Value = "Move Robot42 to location (9, 1) and remove the large debris."
$ws.Range("E3").Value = "(9, 1)"

# Row 4: Robot48/broom -> Robot50/vacuum, location (5, 6) -> (7, 11)
$ws.Range("A4").Value = "Move Robot50 to location (7, 11) and remove the dust."
$ws.Range("B4").Value = "['Robot50']"
$ws.Range("C4").Value = "['vacuum']"
$ws.Range("E4").Value = "(7, 11)"

# Row 5: Robot31 -> Robot41, location (9, 4) -> (6, 12)
$ws.Range("A5").Value = "Move Robot41 to location (6, 12) and remove the grass."
$ws.Range("B5").Value = "['Robot41']"
$ws.Range("E5").Value = "(6, 12)"

# Row 6: Robot8/broom -> Robot50/vacuum, location (8, 12) -> (3, 1)
$ws.Range("A6").Value = "Move Robot50 to location (3, 1) and remove the small debris."
$ws.Range("B6").Value = "['Robot50']"
$ws.Range("C6").Value = "['vacuum']"
$ws.Range("E6").Value = "(3, 1)"

# Row 7: Robot23 -> Robot13, location (11, 1) -> (1, 4)
$ws.Range("A7").Value = "Move Robot13 to location (1, 4) and remove the vehicle."
$ws.Range("B7").Value = "['Robot13']"
$ws.Range("E7").Value = "(1, 4)"

# Row 8: Robot23 -> Robot13, location (12, 10) -> (11, 1)
$ws.Range("A8").Value = "Move Robot13 to location (11, 1) and remove the construction materials."
$ws.Range("B8").Value = "['Robot13']"
$ws.Range("E8").Value = "(11, 1)"

# Row 9: Robot14 location (7, 11) -> (2, 10)
$ws.Range("A9").Value = "Move Robot14 to location (2, 10) and remove the tree branches."
$ws.Range("E9").Value = "(2, 10)"

# Row 10: Robot15 location (5, 3) -> (8, 6)
$ws.Range("A10").Value = "Move Robot15 to location (8, 6) and remove the screws."
$ws.Range("E10").Value = "(8, 6)"
